$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the content of a [start,end) Range with a literal run of
# OOXML (runs, proofErr markers, bookmarks, etc.), by wrapping it in a
# single-part WordprocessingML package and calling Range.InsertXML on a
# freshly-fetched Range object (re-fetching avoids a stale-range quirk when
# the same Range object was just used for Find.Execute).
# ---------------------------------------------------------------------------
function Set-RangeXml($doc, $startPos, $endPos, $xmlBody) {
    $fresh = $doc.Range($startPos, $endPos)
    $wrapper = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $xmlBody + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $fresh.InsertXML($wrapper)
}

# ===========================================================================
# 1) Remove the whole "Progress" section: the "Progress" heading paragraph
#    through the "A login page, ..." paragraph (13 paragraphs). The empty
#    paragraph before "Progress" and the empty paragraph after stay, so two
#    blank paragraphs end up adjacent.
# ===========================================================================
$pStart = $d.Paragraphs(8)
$pEnd = $d.Paragraphs(20)
$delRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$delRange.Delete()

# ===========================================================================
# 2) Rewrite the "... test'. This will run ..." sentence, splitting it into
#    several runs and marking "non UI" with a gramStart/gramEnd proofErr
#    pair, matching the target markup exactly.
# ===========================================================================
$rng = $d.Content
$oldText = " test" + [char]0x2019 + ". This will run both the UI acceptance tests and the acceptance tests that only interact with the system" + [char]0x2019 + "s controllers and models."
$found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$newRuns = ''
$newRuns += '<w:r><w:t xml:space="preserve"> test' + [char]0x2019 + '. This will run both the UI acceptance tests </w:t></w:r>'
$newRuns += '<w:r><w:t>(</w:t></w:r>'
$newRuns += '<w:r><w:t>that</w:t></w:r>'
$newRuns += '<w:r><w:t xml:space="preserve"> test everything)</w:t></w:r>'
$newRuns += '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
$newRuns += '<w:r><w:t xml:space="preserve">and the </w:t></w:r>'
$newRuns += '<w:proofErr w:type="gramStart"/>'
$newRuns += '<w:r><w:t>non UI</w:t></w:r>'
$newRuns += '<w:proofErr w:type="gramEnd"/>'
$newRuns += '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
$newRuns += '<w:r><w:t>acceptance tests that only interact with the system' + [char]0x2019 + 's models.</w:t></w:r>'

Set-RangeXml $d $rng.Start $rng.End $newRuns

# ===========================================================================
# 3) Delete "The non-UI Acceptance Tests cover ..." paragraph plus the empty
#    paragraph that follows it.
# ===========================================================================
$rng2 = $d.Content
$rng2.Find.Execute("The non-UI Acceptance Tests cover as much of the system" + [char]0x2019 + "s specified functionality without user input.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$nonUiPara = $rng2.Paragraphs(1)
$emptyPara = $nonUiPara.Next()
$delRange2 = $d.Range($nonUiPara.Range.Start, $emptyPara.Range.End)
$delRange2.Delete()

# ===========================================================================
# 4) Move the "_GoBack" bookmark to the very start of "The UI Acceptance
#    Tests add cameras ..." paragraph (it previously sat in the
#    now-deleted "We have a UI ..." paragraph).
# ===========================================================================
$rng3 = $d.Content
$rng3.Find.Execute("The UI Acceptance Tests add cameras", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$uiTestsStart = $rng3.Start
$bmRange = $d.Range($uiTestsStart, $uiTestsStart)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ===========================================================================
# 5) Move <w:lastRenderedPageBreak/> from the run starting "This system was
#    constructed ..." to the run starting "The Client is a ".
# ===========================================================================
$rng4 = $d.Content
$rng4.Find.Execute("This system was constructed with the ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Set-RangeXml $d $rng4.Start $rng4.End '<w:r><w:t xml:space="preserve">This system was constructed with the </w:t></w:r>'

$rng5 = $d.Content
$rng5.Find.Execute("The Client is a ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Set-RangeXml $d $rng5.Start $rng5.End '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">The Client is a </w:t></w:r>'

Write-Host "Edit complete."
